$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header column, formatted like the other header cells (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save values: 1 for row 6 (2024-07-14), 0 for every other data row (2..17)
for ($r = 2; $r -le 17; $r++) {
    if ($r -eq 6) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
